$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing XPath locator values (sharedStrings text edits) ---
# EmailAddress locator: placeholder-based -> id-based
$ws.Range("B2").Value = "//input[@id='login_username']"
# password locator: placeholder-based -> id-based
$ws.Range("B3").Value = "//input[@id='login_password']"
# login_Button locator: Login -> Sign in
$ws.Range("B4").Value = "//button[normalize-space()='Sign in']"

# --- Add a new logout_Button key/xpath pair into the now-available row 12 ---
$rowStyle = $ws.Range("A11").Style
$ws.Range("A12").Value = "logout_Button"
$ws.Range("B12").Value = "//span[normalize-space()='logout']"
$ws.Range("A12").Style = $rowStyle
$ws.Range("B12").Style = $rowStyle

# --- Move the active selection to B13 ---
[void]$ws.Range("B13").Select()

# --- Remove the now-unused trailing blank rows 13:16 ---
$ws.Range("A13:B16").EntireRow.Delete()
